$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily record was added at the top of the table (row 2). This pushes
# every existing data row down by one. The sheet's used range grows from
# A1:T38 to A1:T39 (the last existing row, old row 38, ends up at row 39).
#
# Shift rows 38..2 down into rows 39..3, working from the bottom up so we
# never clobber a row before it has been copied.
for ($r = 38; $r -ge 2; $r--) {
    $nr = $r + 1
    for ($c = 1; $c -le 20; $c++) {
        $src = $ws.Cells.Item($r, $c)
        $dst = $ws.Cells.Item($nr, $c)
        $dst.Value = $src.Value2
    }
}

# The brand-new last row (39) needs the same date format as every other row
# in column D (style carries numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Cells.Item(39, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Now populate the new row 2 with the freshly reported data. Most columns
# (market/product metadata) are identical to what used to be in row 2;
# only the date and the price-related columns differ.
$ws.Cells.Item(2, 1).Value = 4
$ws.Cells.Item(2, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(2, 3).Value = "Los Lagos"
$ws.Cells.Item(2, 4).Value = 44922
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 5).Value = 10
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100101
$ws.Cells.Item(2, 8).Value = "Berries"
$ws.Cells.Item(2, 9).Value = 100101001
$ws.Cells.Item(2, 10).Value = "Arándano (blue)"
$ws.Cells.Item(2, 11).Value = "Sin especificar"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 400
$ws.Cells.Item(2, 14).Value = 3500
$ws.Cells.Item(2, 15).Value = 3800
$ws.Cells.Item(2, 16).Value = 3650
$ws.Cells.Item(2, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(2, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(2, 19).Value = 1825
$ws.Cells.Item(2, 20).Value = 2
